$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C2 to a numeric value (was an empty inline string cell)
$ws.Range("C2").Value = 2601719.61

# Row 5: reset C5:F5 to 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Row 6: reset C6:F6 to 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
